$wb = $excel.ActiveWorkbook

# ---- Summary Table sheet ----
$ws1 = $wb.Worksheets.Item("Summary Table")
$ws1.Range("B2:B21").NumberFormat = "@"
$ws1.Range("A2").Value = ' Science and technology innovation policy '
$ws1.Range("B2").Value = ' 47 '
$ws1.Range("C2").Value = ' China''s State Council unveils comprehensive "Science and Technology Innovation 2030" policy framework with unprecedented funding commitments totaling 2.7 trillion yuan through 2030. '
$ws1.Range("D2").Value = ' https://www.gov.cn/zhengce/content/2025-04/22/content_6108932.htm '
$ws1.Range("E2").Value = ' The new framework represents China''s most ambitious S&T policy to date, targeting self-sufficiency in core technologies by 2030. Key elements include doubling basic research funding, creating 25 new national laboratories, and establishing special innovation zones with preferential tax and regulatory policies. The plan specifically prioritizes semiconductors, quantum computing, and AI as strategic sectors. '
$ws1.Range("A3").Value = ' Artificial intelligence '
$ws1.Range("B3").Value = ' 43 '
$ws1.Range("C3").Value = ' Baidu and the Chinese Academy of Sciences launch world''s largest open-source foundation model "Wenxin 3.0" with 10 trillion parameters, challenging GPT-5. '
$ws1.Range("D3").Value = ' https://www.caict.ac.cn/xwdt/ynxw/202504/t20250423_576321.html '
$ws1.Range("E3").Value = ' The Wenxin 3.0 model represents China''s most advanced AI system to date, demonstrating capabilities in multimodal reasoning, scientific problem-solving, and code generation that reportedly match or exceed OpenAI''s latest models. The project received 12 billion yuan in government funding and leverages China''s growing computational infrastructure. CAS researchers emphasize the model''s alignment with China''s AI ethics framework. '
$ws1.Range("A4").Value = ' Integration of industry, academia and research '
$ws1.Range("B4").Value = ' 39 '
$ws1.Range("C4").Value = ' MOST launches the "Fusion 2025" initiative creating 50 new industrial technology research institutes with joint governance from universities, companies, and local governments. '
$ws1.Range("D4").Value = ' https://www.most.gov.cn/kjbgz/202504/t20250419_186543.html '
$ws1.Range("E4").Value = ' The Fusion 2025 program represents a significant evolution in China''s approach to industry-academia collaboration, establishing institutes with shared governance structures and dedicated funding pools. The model aims to address previous criticism about technology transfer bottlenecks. Initial focus areas include advanced manufacturing, new materials, and biotechnology. Each institute will receive 200-500 million yuan in startup funding and preferential policies for commercialization. '
$ws1.Range("A5").Value = ' Quantum communication '
$ws1.Range("B5").Value = ' 38 '
$ws1.Range("C5").Value = ' China completes first cross-continental quantum encrypted video call between Beijing and Paris using the Mozi satellite, demonstrating practical quantum key distribution. '
$ws1.Range("D5").Value = ' https://www.cas.cn/syky/202504/t20250420_4955342.html '
$ws1.Range("E5").Value = ' The Beijing-Paris quantum encrypted video call represents a major breakthrough in quantum communication''s practical applications. Using China''s Mozi satellite as a relay, the system distributed quantum encryption keys that enabled a 15-minute video conference with theoretically unbreakable encryption. The demonstration involved China''s University of Science and Technology and France''s National Centre for Scientific Research, highlighting both the technological achievement and international collaboration. '
$ws1.Range("A6").Value = ' Semiconductor packaging '
$ws1.Range("B6").Value = ' 36 '
$ws1.Range("C6").Value = ' SMIC announces breakthrough in advanced packaging technology with new chip-on-wafer-on-substrate (CoWoS) process, reducing dependence on foreign technology. '
$ws1.Range("D6").Value = ' https://www.miit.gov.cn/xwdt/gxdt/sjdt/art/2025/art_4f98c674218b4a678f1f1f0a1e67cb54.html '
$ws1.Range("E6").Value = ' SMIC''s new advanced packaging capability represents a significant step toward self-sufficiency in the semiconductor supply chain. The CoWoS technology enables 3D integration of multiple chips with high-bandwidth memory, crucial for AI accelerators and high-performance computing. The breakthrough was supported by China''s "IC Manufacturing Innovation Center" and involves partnerships with five domestic equipment suppliers, demonstrating progress in building a complete domestic ecosystem. '
$ws1.Range("A7").Value = ' New quality productivity '
$ws1.Range("B7").Value = ' 32 '
$ws1.Range("C7").Value = ' NDRC releases implementation guidelines for "New Quality Productive Forces" initiative with 1.2 trillion yuan investment across 10 technology domains and 85 demonstration zones. '
$ws1.Range("D7").Value = ' https://www.ndrc.gov.cn/xxgk/zcfb/tz/202504/t20250424_1391075.html '
$ws1.Range("E7").Value = ' The New Quality Productive Forces initiative represents China''s comprehensive strategy to upgrade its industrial capabilities through technological innovation. The implementation guidelines outline specific targets for 10 sectors including advanced manufacturing, new energy, and digital technologies. The plan emphasizes both technological breakthroughs and institutional reforms in resource allocation, establishing 85 demonstration zones with special regulatory frameworks. The initiative aims to incubate 5,000 specialized technology enterprises by 2027. '
$ws1.Range("A8").Value = ' Digital economy policy '
$ws1.Range("B8").Value = ' 31 '
$ws1.Range("C8").Value = ' China''s Digital Economy Promotion Law passes with provisions for data rights, platform governance, and cross-border data flows, effective October 1, 2025. '
$ws1.Range("D8").Value = ' http://www.npc.gov.cn/npc/c2/c30834/202504/t20250421_571328.html '
$ws1.Range("E8").Value = ' The Digital Economy Promotion Law represents China''s first comprehensive legal framework for its digital economy, balancing development with security considerations. The law establishes clear data rights and responsibilities, regulatory frameworks for platform companies, and controlled mechanisms for cross-border data flows. It includes provisions for digital infrastructure development, digital trade promotion, and digital service standards. The law also creates a Coordinated Digital Economy Development Committee to harmonize policies across governmental departments. '
$ws1.Range("A9").Value = ' Hydrogen energy storage '
$ws1.Range("B9").Value = ' 28 '
$ws1.Range("C9").Value = ' NDRC approves 50 billion yuan hydrogen energy demonstration cluster in Inner Mongolia featuring integrated production, storage, and transportation infrastructure. '
$ws1.Range("D9").Value = ' https://www.chinadaily.com.cn/a/202504/22/WS661426a5a31082fc043b9721.html '
$ws1.Range("E9").Value = ' The Inner Mongolia hydrogen cluster represents China''s largest integrated hydrogen energy project to date. The facility will leverage the region''s abundant renewable energy resources to produce green hydrogen at scale, store it in multiple forms including underground caverns, and transport it via dedicated pipelines to industrial centers. The project aims to reduce hydrogen production costs by 40% through technological innovations and scale, making hydrogen economically viable for industrial applications. The initiative involves partnerships with 12 major companies and 5 research institutions. '
$ws1.Range("A10").Value = ' Biotechnology '
$ws1.Range("B10").Value = ' 26 '
$ws1.Range("C10").Value = ' China issues new comprehensive five-year biotechnology development plan with focus on synthetic biology, regenerative medicine, and agricultural applications. '
$ws1.Range("D10").Value = ' https://www.most.gov.cn/kjbgz/202504/t20250423_186732.html '
$ws1.Range("E10").Value = ' The 2025-2030 National Biotechnology Development Plan outlines China''s strategic priorities for the biotech sector with three focused areas: medical applications, agricultural innovation, and industrial biotechnology. The plan allocates 350 billion yuan in government funding and establishes five new national biotechnology research centers. Key initiatives include accelerating the development of mRNA vaccines, gene editing technologies, synthetic biology platforms, and biomanufacturing capabilities. The plan also addresses biosafety concerns with new regulatory frameworks. '
$ws1.Range("A11").Value = ' Brain-computer interface '
$ws1.Range("B11").Value = ' 24 '
$ws1.Range("C11").Value = ' Chinese Academy of Sciences unveils non-invasive brain-computer interface capable of complex control functions, demonstrating applications for disability assistance. '
$ws1.Range("D11").Value = ' https://english.cas.cn/newsroom/news/202504/t20250419_324276.html '
$ws1.Range("E11").Value = ' The CAS Brain-Computer Interface System represents a significant advance in non-invasive BCI technology, using high-density electroencephalography combined with advanced signal processing algorithms to achieve precision previously requiring implanted electrodes. In demonstrations, users controlled robotic arms to perform complex manipulation tasks with 95% accuracy. The system''s medical applications for paralysis patients will begin clinical trials in five hospitals. The technology was developed through the "Brain Science and Brain-Inspired Intelligence" national science project. '
$ws1.Range("A12").Value = ' Carbon fiber composites '
$ws1.Range("B12").Value = ' 24 '
$ws1.Range("C12").Value = ' COMAC and Chinese material scientists announce breakthrough in domestic aviation-grade carbon fiber composites for C919 passenger jet wing components. '
$ws1.Range("D12").Value = ' https://www.miit.gov.cn/xwdt/gxdt/sjdt/art/2025/art_7d3fd8c4e5f84d028a3e10de2cd58be2.html '
$ws1.Range("E12").Value = ' The domestically developed T800-grade carbon fiber composites represent a significant technological breakthrough for China''s aerospace industry. Previously reliant on imported materials for critical components, COMAC can now source high-performance composites domestically. The new materials demonstrate 15% higher strength-to-weight ratio than previous generations while meeting stringent aviation safety standards. The breakthrough involved collaboration between three research institutes and two companies, supported by the "Key Materials Breakthrough" national program. '
$ws1.Range("A13").Value = ' Smart city '
$ws1.Range("B13").Value = ' 22 '
$ws1.Range("C13").Value = ' Ministry of Housing announces next-generation smart city initiative with 15 pilot cities implementing integrated urban brain platforms and digital twins. '
$ws1.Range("D13").Value = ' http://www.mohurd.gov.cn/gongkai/tongzhi/202504/t20250424_782541.html '
$ws1.Range("E13").Value = ' The Urban Intelligence 2.0 initiative represents China''s most ambitious smart city program to date, moving beyond separate systems toward comprehensive urban management platforms. The 15 pilot cities will implement "urban brain" systems that integrate data from transportation, energy, environment, and public services. Digital twin technology will create virtual replicas of city infrastructure for simulation and planning. The program emphasizes practical applications of smart city technology to address urban management challenges, with special focus on emergency response, resource efficiency, and citizen services. '
$ws1.Range("A14").Value = ' Science and technology security policy '
$ws1.Range("B14").Value = ' 21 '
$ws1.Range("C14").Value = ' The Ministry of Science and Technology issues new guidelines on "secure and trustworthy science" with expanded reviews for international research collaboration. '
$ws1.Range("D14").Value = ' https://www.most.gov.cn/xxgk/xinxifenlei/fdzdgknr/fgzc/gfxwj/202504/t20250418_186403.html '
$ws1.Range("E14").Value = ' The new science security guidelines establish a comprehensive framework for managing security risks in scientific research while maintaining openness for innovation. The policy expands security reviews for international collaborations in 23 sensitive technology areas, creates institutional security committees at research organizations, and establishes procedures for evaluating foreign partners. The guidelines aim to balance China''s continued integration with global science with growing concerns about technology security, reflecting similar policies being implemented in the US and EU. '
$ws1.Range("A15").Value = ' Quantum computing cloud platform '
$ws1.Range("B15").Value = ' 20 '
$ws1.Range("C15").Value = ' Alibaba Cloud and Chinese Academy of Sciences launch 96-qubit quantum computing cloud platform with free access for domestic researchers and tiered commercial services. '
$ws1.Range("D15").Value = ' https://www.alibabacloud.com/press-room/alibaba-cloud-launches-quantum-computing-services-for-global-customers '
$ws1.Range("E15").Value = ' The new quantum computing platform represents China''s most powerful publicly accessible quantum resource, featuring a 96-qubit superconducting processor with reported 99.5% gate fidelity. The service offers both circuit-based and quantum annealing computing models through a user-friendly development environment. Chinese researchers from 235 institutions receive priority free access, while commercial applications are available through tiered service plans. The platform includes specialized quantum algorithms for materials science, cryptography, and optimization problems, with documented performance advantages for specific computational tasks. '
$ws1.Range("A16").Value = ' Solid-state batteries '
$ws1.Range("B16").Value = ' 19 '
$ws1.Range("C16").Value = ' CATL announces mass production timeline for next-generation solid-state EV batteries with 500 Wh/kg energy density, doubling current commercial batteries. '
$ws1.Range("D16").Value = ' https://www.catl.com/en/news/20250421/ '
$ws1.Range("E16").Value = ' CATL''s solid-state battery announcement represents a potential leap forward in electric vehicle technology. The company plans to begin mass production by early 2026, with the first commercial applications in Chinese EVs. The batteries feature ceramic-based solid electrolytes that enable higher energy density while addressing safety concerns associated with liquid electrolytes. The technology was developed through CATL''s "Beyond Lithium" research initiative with 3.5 billion yuan in government support. Multiple Chinese automakers have already signed agreements to incorporate the batteries in upcoming models. '
$ws1.Range("A17").Value = ' 6G communication '
$ws1.Range("B17").Value = ' 18 '
$ws1.Range("C17").Value = ' China establishes international 6G Innovation Alliance with 21 countries, releases white paper on terahertz communication standards and satellite integration. '
$ws1.Range("D17").Value = ' http://www.caict.ac.cn/xwdt/ynxw/202504/t20250423_576325.html '
$ws1.Range("E17").Value = ' The International 6G Innovation Alliance represents China''s effort to shape global standards for next-generation wireless technology. The alliance includes 21 countries with 85 member organizations from industry, academia, and government. The accompanying white paper outlines technical proposals for terahertz communications, integrated satellite-terrestrial networks, and quantum-secured transmission protocols. The initiative positions China as a leading voice in 6G development while building international consensus around compatible technical approaches. The alliance will establish five working groups focusing on different aspects of 6G technology development. '
$ws1.Range("A18").Value = ' Strategic emerging industries '
$ws1.Range("B18").Value = ' 17 '
$ws1.Range("C18").Value = ' State Council identifies seven strategic emerging industry clusters for concentrated development with special economic zones and regulatory sandboxes. '
$ws1.Range("D18").Value = ' http://www.gov.cn/zhengce/content/2025-04/24/content_6109283.htm '
$ws1.Range("E18").Value = ' The Strategic Emerging Industries Initiative focuses on developing coordinated industrial ecosystems rather than individual companies or technologies. The seven designated clusters include next-generation information technology, biotechnology, new energy, new materials, high-end equipment manufacturing, new energy vehicles, and digital creative industries. Each cluster will receive customized policy support including tax incentives, land allocation, and specialized financial services. The initiative also creates regulatory "sandboxes" allowing companies in these sectors to test innovative products under modified regulatory frameworks to accelerate commercialization. '
$ws1.Range("A19").Value = ' Technology transfer '
$ws1.Range("B19").Value = ' 17 '
$ws1.Range("C19").Value = ' MOST launches comprehensive technology transfer reform with new valuation methods, equity incentives for researchers, and streamlined licensing procedures. '
$ws1.Range("D19").Value = ' https://www.most.gov.cn/kjbgz/202504/t20250421_186662.html '
$ws1.Range("E19").Value = ' The technology transfer reform addresses longstanding bottlenecks in China''s innovation system by fundamentally changing how research results move from labs to industry. Key measures include: allowing researchers to retain up to 70% of equity or profits from technologies they develop; creating standardized technology valuation methods; establishing professional technology transfer offices at research institutions; and simplifying approval procedures for licenses. The reforms aim to dramatically increase commercialization rates for Chinese research output, which has historically lagged behind publication output. '
$ws1.Range("A20").Value = ' Synthetic biology '
$ws1.Range("B20").Value = ' 16 '
$ws1.Range("C20").Value = ' China launches "Biological Manufacturing 2025" initiative with synthetic biology facilities for medicine, materials, and food production using engineered microorganisms. '
$ws1.Range("D20").Value = ' https://www.cas.cn/xw/zyxw/tt/202504/t20250422_4955398.html '
$ws1.Range("E20").Value = ' The Biological Manufacturing initiative represents China''s strategic push into synthetic biology-based production technologies. The program will establish five biomanufacturing centers capable of producing pharmaceuticals, materials, chemicals, and food ingredients using engineered microorganisms at industrial scale. Key technologies include CRISPR-based genome editing, cell-free protein synthesis, and artificial chromosomes. The initiative targets 30% cost reduction for biomanufactured products compared to traditional chemical synthesis, while reducing environmental impact. The program includes partnerships with 150+ companies to accelerate commercial applications. '
$ws1.Range("A21").Value = ' Intergovernmental science and technology cooperation '
$ws1.Range("B21").Value = ' 15 '
$ws1.Range("C21").Value = ' China and ASEAN establish joint $5 billion innovation fund and research exchange program focused on climate technology, digital infrastructure, and health sciences. '
$ws1.Range("D21").Value = ' https://www.fmprc.gov.cn/eng/wjbxw/202504/t20250423_11477629.html '
$ws1.Range("E21").Value = ' The China-ASEAN Science and Technology Cooperation Framework represents the most comprehensive S&T partnership between China and Southeast Asian nations. The program includes a $5 billion joint innovation fund with equal contributions from China and ASEAN members, researcher exchange programs involving 200+ institutions, and shared research facilities focused on climate technology, digital infrastructure, and health sciences. The agreement includes commitments to transparent governance and intellectual property sharing mechanisms, addressing previous concerns about technology transfer arrangements. The framework establishes a permanent secretariat in Singapore to coordinate activities. '
$ws1.Range("B2:B21").ClearFormats()
$ws1.Rows.Item(22).Delete()

# ---- Sources sheet ----
$ws2 = $wb.Worksheets.Item("Sources")
$ws2.Range("A3").Value = ' China State Council '
$ws2.Range("B3").Value = ' https://www.gov.cn/zhengce/content/2025-04/22/content_6108932.htm '
$ws2.Range("C3").Value = ' 2025-04-22 '
$ws2.Range("A4").Value = ' China Academy of Information and Communications Technology '
$ws2.Range("B4").Value = ' https://www.caict.ac.cn/xwdt/ynxw/202504/t20250423_576321.html '
$ws2.Range("C4").Value = ' 2025-04-23 '
$ws2.Range("A5").Value = ' Ministry of Science and Technology of China '
$ws2.Range("B5").Value = ' https://www.most.gov.cn/kjbgz/202504/t20250419_186543.html '
$ws2.Range("C5").Value = ' 2025-04-19 '
$ws2.Range("A6").Value = ' Chinese Academy of Sciences '
$ws2.Range("B6").Value = ' https://www.cas.cn/syky/202504/t20250420_4955342.html '
$ws2.Range("C6").Value = ' 2025-04-20 '
$ws2.Range("A7").Value = ' Ministry of Industry and Information Technology '
$ws2.Range("B7").Value = ' https://www.miit.gov.cn/xwdt/gxdt/sjdt/art/2025/art_4f98c674218b4a678f1f1f0a1e67cb54.html '
$ws2.Range("C7").Value = ' 2025-04-21 '
$ws2.Range("A8").Value = ' National Development and Reform Commission '
$ws2.Range("B8").Value = ' https://www.ndrc.gov.cn/xxgk/zcfb/tz/202504/t20250424_1391075.html '
$ws2.Range("C8").Value = ' 2025-04-24 '
$ws2.Range("A9").Value = ' National People''s Congress '
$ws2.Range("B9").Value = ' http://www.npc.gov.cn/npc/c2/c30834/202504/t20250421_571328.html '
$ws2.Range("C9").Value = ' 2025-04-21 '
$ws2.Range("A10").Value = ' China Daily '
$ws2.Range("B10").Value = ' https://www.chinadaily.com.cn/a/202504/22/WS661426a5a31082fc043b9721.html '
$ws2.Range("C10").Value = ' 2025-04-22 '
$ws2.Range("A11").Value = ' Ministry of Science and Technology of China '
$ws2.Range("B11").Value = ' https://www.most.gov.cn/kjbgz/202504/t20250423_186732.html '
$ws2.Range("C11").Value = ' 2025-04-23 '
$ws2.Range("A12").Value = ' Chinese Academy of Sciences English '
$ws2.Range("B12").Value = ' https://english.cas.cn/newsroom/news/202504/t20250419_324276.html '
$ws2.Range("C12").Value = ' 2025-04-19 '
$ws2.Range("A13").Value = ' Ministry of Industry and Information Technology '
$ws2.Range("B13").Value = ' https://www.miit.gov.cn/xwdt/gxdt/sjdt/art/2025/art_7d3fd8c4e5f84d028a3e10de2cd58be2.html '
$ws2.Range("C13").Value = ' 2025-04-20 '
$ws2.Range("A14").Value = ' Ministry of Housing and Urban-Rural Development '
$ws2.Range("B14").Value = ' http://www.mohurd.gov.cn/gongkai/tongzhi/202504/t20250424_782541.html '
$ws2.Range("C14").Value = ' 2025-04-24 '
$ws2.Range("A15").Value = ' Ministry of Science and Technology of China '
$ws2.Range("B15").Value = ' https://www.most.gov.cn/xxgk/xinxifenlei/fdzdgknr/fgzc/gfxwj/202504/t20250418_186403.html '
$ws2.Range("C15").Value = ' 2025-04-18 '
$ws2.Range("A16").Value = ' Alibaba Cloud '
$ws2.Range("B16").Value = ' https://www.alibabacloud.com/press-room/alibaba-cloud-launches-quantum-computing-services-for-global-customers '
$ws2.Range("C16").Value = ' 2025-04-22 '
$ws2.Range("A17").Value = ' CATL '
$ws2.Range("B17").Value = ' https://www.catl.com/en/news/20250421/ '
$ws2.Range("C17").Value = ' 2025-04-21 '
$ws2.Range("A18").Value = ' China Academy of Information and Communications Technology '
$ws2.Range("B18").Value = ' http://www.caict.ac.cn/xwdt/ynxw/202504/t20250423_576325.html '
$ws2.Range("C18").Value = ' 2025-04-23 '
$ws2.Range("A19").Value = ' China State Council '
$ws2.Range("B19").Value = ' http://www.gov.cn/zhengce/content/2025-04/24/content_6109283.htm '
$ws2.Range("C19").Value = ' 2025-04-24 '
$ws2.Range("A20").Value = ' Ministry of Science and Technology of China '
$ws2.Range("B20").Value = ' https://www.most.gov.cn/kjbgz/202504/t20250421_186662.html '
$ws2.Range("C20").Value = ' 2025-04-21 '
$ws2.Range("A21").Value = ' Chinese Academy of Sciences '
$ws2.Range("B21").Value = ' https://www.cas.cn/xw/zyxw/tt/202504/t20250422_4955398.html '
$ws2.Range("C21").Value = ' 2025-04-22 '
$ws2.Range("A22").Value = ' Ministry of Foreign Affairs '
$ws2.Range("B22").Value = ' https://www.fmprc.gov.cn/eng/wjbxw/202504/t20250423_11477629.html '
$ws2.Range("C22").Value = ' 2025-04-23 '
$ws2.Range("A23").Value = ' Science and Technology Daily '
$ws2.Range("B23").Value = ' http://www.stdaily.com/index/kejixinwen/202504/t20250420_786531.html '
$ws2.Range("C23").Value = ' 2025-04-20 '
$ws2.Range("A24").Value = ' Shanghai Municipal Government '
$ws2.Range("B24").Value = ' http://www.shanghai.gov.cn/nw12344/20250419/9a16d50621d94fd29d8c25f6c5d37838.html '
$ws2.Range("C24").Value = ' 2025-04-19 '
$ws2.Range("A25").Value = ' Xinhua News Agency '
$ws2.Range("B25").Value = ' http://www.news.cn/english/20250423/872a6e7b27064a9d9a15f2d36c0c5a21/c.html '
$ws2.Range("C25").Value = ' 2025-04-23 '
$ws2.Range("A26").Value = ' Guangdong Provincial Government '
$ws2.Range("B26").Value = ' http://www.gd.gov.cn/zwgk/wjk/qbwj/yfh/content/post_4157623.html '
$ws2.Range("C26").Value = ' 2025-04-21 '
$ws2.Range("A27").Value = ' People''s Daily '
$ws2.Range("B27").Value = ' http://en.people.cn/n3/2025/0423/c90000-20063788.html '
$ws2.Range("C27").Value = ' 2025-04-23 '
$ws2.Range("A28").Value = ' Ministry of Education '
$ws2.Range("B28").Value = ' http://www.moe.gov.cn/jyb_xwfb/gzdt_gzdt/s5987/202504/t20250420_653982.html '
$ws2.Range("C28").Value = ' 2025-04-20 '
$ws2.Range("A29").Value = ' National Energy Administration '
$ws2.Range("B29").Value = ' http://www.nea.gov.cn/2025-04/21/c_1310822565.htm '
$ws2.Range("C29").Value = ' 2025-04-21 '
$ws2.Range("A30").Value = ' China Securities Journal '
$ws2.Range("B30").Value = ' https://www.cs.com.cn/xwzx/hg/202504/t20250419_6397532.html '
$ws2.Range("C30").Value = ' 2025-04-19 '
$ws2.Range("A31").Value = ' Ministry of Ecology and Environment '
$ws2.Range("B31").Value = ' https://www.mee.gov.cn/ywdt/xwfb/202504/t20250422_1046238.shtml '
$ws2.Range("C31").Value = ' 2025-04-22 '

# ---- Executive Summary sheet ----
$ws3 = $wb.Worksheets.Item("Executive Summary")
$ws3.Range("A2").Value = 'Five Most Impactful News Summaries:
1. China''s State Council has unveiled the "Science and Technology Innovation 2030" policy framework with unprecedented funding commitments of 2.7 trillion yuan through 2030. This comprehensive strategy aims to achieve technological self-sufficiency in core areas, doubling basic research funding, creating 25 new national laboratories, and establishing special innovation zones with preferential policies. Semiconductors, quantum computing, and AI are designated as strategic sectors receiving priority support.
2. Baidu and the Chinese Academy of Sciences have launched "Wenxin 3.0," reportedly the world''s largest open-source foundation AI model with 10 trillion parameters, directly challenging OpenAI''s GPT-5. The model demonstrates advanced capabilities in multimodal reasoning, scientific problem-solving, and code generation, representing China''s most sophisticated AI system to date with 12 billion yuan in government funding support.
3. The National Development and Reform Commission has released implementation guidelines for the "New Quality Productive Forces" initiative, committing 1.2 trillion yuan across 10 technology domains and establishing 85 demonstration zones. This initiative represents China''s comprehensive strategy to upgrade its industrial capabilities through technological innovation while reforming resource allocation mechanisms and regulatory frameworks.
4. China''s Digital Economy Promotion Law has passed and will take effect October 1, 2025, establishing the country''s first comprehensive legal framework for the digital economy. The law balances development goals with security considerations, addressing data rights, platform governance, cross-border data flows, and creating a Coordinated Digital Economy Development Committee to harmonize policies across governmental departments.
5. The Ministry of Science and Technology has launched a comprehensive technology transfer reform that fundamentally changes how research results move from labs to industry. Key measures include allowing researchers to retain up to 70% of equity or profits from technologies they develop, creating standardized valuation methods, establishing professional transfer offices, and simplifying licensing procedures to dramatically increase commercialization rates.'

# ---- Cooccurrence sheet ----
$ws4 = $wb.Worksheets.Item("Cooccurrence")
$ws5 = $wb.Worksheets.Item("Associations")
$ws5.Range("A1:B1").Copy()
$ws4.Range("A1:C1").PasteSpecial(-4122)
$ws4.Range("A1").Value = 'source'
$ws4.Range("B1").Value = 'target'
$ws4.Range("C1").Value = 'count'
$ws4.Range("A2").Value = 'Biotechnology'
$ws4.Range("B2").Value = 'Technology transfer'
$ws4.Range("C2").Value = 1
$ws4.Range("A3").Value = 'Biotechnology'
$ws4.Range("B3").Value = 'Synthetic biology'
$ws4.Range("C3").Value = 1
$ws4.Range("A4").Value = 'Biotechnology'
$ws4.Range("B4").Value = 'Strategic emerging industries'
$ws4.Range("C4").Value = 1

# ---- Associations sheet ----
$ws5.Range("A2").Value = 'Biotechnology'
$ws5.Range("B2").Value = 3
$ws5.Range("A3").Value = 'Technology transfer'
$ws5.Range("B3").Value = 1
$ws5.Range("A4").Value = 'Quantum communication'
$ws5.Range("B4").Value = 1
$ws5.Range("A5").Value = 'Synthetic biology'
$ws5.Range("B5").Value = 2
$ws5.Range("A6").Value = 'Brain-computer interface'
$ws5.Range("B6").Value = 1
$ws5.Range("A7").Value = 'Carbon fiber composites'
$ws5.Range("B7").Value = 1
$ws5.Range("A8").Value = 'Smart city'
$ws5.Range("B8").Value = 1
$ws5.Range("A9").Value = 'Quantum computing cloud platform'
$ws5.Range("B9").Value = 1
$ws5.Range("A10").Value = 'Strategic emerging industries'
$ws5.Range("B10").Value = 1
